# Auto-generated COM-interop script to rewrite the mock-encounter
# transcripts (enc_1/enc_2/enc_3) and remove the now-unused enc_4 sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('enc_1')
$ws.Cells.Item(2, 1).Value = 'Hi Gabriella, I understand that you''re here today because you''ve been experiencing some low mood.'
$ws.Cells.Item(2, 2).Value = 'Yeah, I just recently I''ve just... I''m just feeling down a lot and I can''t really, like, snap out of it. Like, sometimes if it''s happened before, I can just, like, get myself back on track. But recently I''m just finding it really difficult to do.'
$ws.Cells.Item(3, 1).Value = 'So it sounds like you''ve been experiencing some low mood and you''re kind of struggling to get going a bit.'
$ws.Cells.Item(3, 2).Value = 'Yeah, definitely.'
$ws.Cells.Item(4, 1).Value = 'You''ve had low mood in the past, but it''s kind of... you''ve been able to pull yourself out of it before, but more recently you''ve found that harder. Okay, so I understand that you first went to your GP, who referred you here about two months ago. Can I ask what led you to go to the GP in the first place?'
$ws.Cells.Item(4, 2).Value = 'My long-term boyfriend broke up with me and it was, it was such a shock, like I really wasn''t expecting it. It just got me really down. I felt really lonely. And I didn''t really want to see people. I just wanted to be in my own space. My mom told me I should probably go to the doctors and he''s referred me to you.'
$ws.Cells.Item(5, 1).Value = 'Okay? So it sounds like you had that breakup and you weren''t expecting it, so it kind of came out of the blue a bit and you''ve been on your own quite a lot. So you''re not been seeing people as much. And it sounds like your mum noticed that there''s been a bit of a change in you. Okay, and how have things been then since then, since that happened and since you went to the GP?'
$ws.Cells.Item(5, 2).Value = 'I just, nothing''s really changed. Nothing''s really improved. I just... I feel down, like everything''s an effort and feel distant from things. And like I said, I can''t snap out of it.'
$ws.Cells.Item(6, 1).Value = 'Okay, and how did you feel about coming here today?'
$ws.Cells.Item(6, 2).Value = 'I was kind of disappointed in myself because I feel like nobody else my age is doing this. I should have been able to snap out of it, like it''s just another thing that I''ve not really... you know, things that''s making me stand out for the wrong reasons.'
$ws.Cells.Item(7, 1).Value = 'Okay, so you feel like it''s kind of odd or unusual to be coming to see a clinical psychologist. I suppose I just want to reassure you that it''s absolutely normal to have feelings of depression and anxiety. It''s really very common, and it certainly isn''t odd or strange at all. So it sounds like it was actually quite hard for you to come here today. So what I guess led you to do that? Because it sounds like maybe there was a lot pushing you against coming here today, but you''ve managed to come anyway.'
$ws.Cells.Item(7, 2).Value = 'I''m just hoping that it''s going to help. I think I didn''t really want to admit it because, like, it''s just another failure, but my mum was talking to me and she''s not really happy with the way I''m behaving. I didn''t really want to disappoint her again. And so I thought there''s no harm in it. So I came.'
$ws.Cells.Item(8, 1).Value = 'So it sounds like partly you''re doing it then for your mum because you''re a bit worried that she isn''t pleased with you at the moment.'
$ws.Cells.Item(8, 2).Value = 'Yeah.'
$ws.Cells.Item(9, 1).Value = 'Well, I think it''s great that you''ve come here. And definitely there''s a lot of evidence to suggest that psychological therapy can help depression, and particularly the type of therapy that we''re gonna be doing, which is called cognitive behavior therapy. I''ll tell you a bit more about that as we go on, but I''m just trying to find out a bit more about you today. So I saw from your notes that you''re a medic student.'
$ws.Cells.Item(9, 2).Value = 'Yeah? I''m not currently in uni at the moment, but yeah, last year I was a medical student.'
$ws.Cells.Item(10, 1).Value = 'Okay, so I guess why, what''s the reason that you''re not in uni this year?'
$ws.Cells.Item(10, 2).Value = 'It was pretty tough first year, so they then suggested that I took a year out. This is the break.'
$ws.Cells.Item(11, 1).Value = 'Okay, so just thinking, so that would be about 6 months ago, was it that your year out started? And it sounds then maybe like your first year was also quite tough?'
$ws.Cells.Item(11, 2).Value = 'Yeah, it was. I thought it''s tough for everyone, but obviously we all take it in different ways, don''t we.'
$ws.Cells.Item(12, 1).Value = 'And you''re planning on returning?'
$ws.Cells.Item(12, 2).Value = 'Yeah, I''m going to have to really.'
$ws.Cells.Item(13, 1).Value = 'But you don''t sound too sure about that.'
$ws.Cells.Item(13, 2).Value = 'No, though I think returning is probably... that''s what''s going to make everyone happy. And I think, you know, at least I''ll be able to finally achieve something and get something right. So yeah, like next year I''m going to... I''ll get through that. Yeah.'
$ws.Cells.Item(14, 1).Value = 'Okay, so it doesn''t sound like you want to go back.'
$ws.Cells.Item(14, 2).Value = 'I do. I didn''t enjoy it. I just... I just don''t feel like me at the moment. And I just, it''s a lot of work. It''s a lot of pressure. But that''s what comes with it. And, you know, it''ll work out for the best if I go.'
$ws.Cells.Item(15, 1).Value = 'Okay, so I guess then there''s quite a lot of thoughts and feelings about going back at the moment. I''m just wondering, in the meantime, what are you up to, what have you been doing?'
$ws.Cells.Item(15, 2).Value = 'Well, I haven''t really had much on the horizon. My mom, her friends own a restaurant and she suggested that I did a bit of waitressing there, just doing it for a bit of money and stuff because obviously I''m not in uni, so you know, I''ve done that.'
$ws.Cells.Item(16, 1).Value = 'Okay, how many shifts have you done there?'
$ws.Cells.Item(16, 2).Value = 'I''ve only done one so far. Yeah, I don''t know whether it''s something that I''d probably go back to.'
$ws.Cells.Item(17, 1).Value = 'Okay, when was that shift?'
$ws.Cells.Item(17, 2).Value = 'It was last week, actually.'
$ws.Cells.Item(18, 1).Value = 'Okay, and how was it? I suppose I''m intrigued, can you tell me a bit more about that?'
$ws.Cells.Item(18, 2).Value = 'It was, it was different, you know, because it was like... I did enjoy it because it''s just a different aspect and I had to speak to people and try to get their orders and try to remember it and stuff. I wasn''t very good, but I suppose it was handy to have some money, like get out of the house.'
$ws.Cells.Item(19, 1).Value = 'So it sounds like you quite enjoyed that. Then what makes you think you weren''t very good?'
$ws.Cells.Item(19, 2).Value = 'I just felt slow. And people, sometimes they didn''t understand what I was saying, so I''d have to talk louder or quieter, slower. And I got a bit confused sometimes with the orders, and I just don''t really think that I''m not really like a friendly... like people don''t really think that I''m a friendly person, so they didn''t really want to talk to me and stuff. I felt like I was inconveniencing them, by like asking them things and getting things wrong. Yeah, they said that apparently most people get like a few mistakes on their first shift, but then again, she''s my mom''s friend, so that''s probably why she''s saying it.'
$ws.Cells.Item(20, 1).Value = 'So the person who owns the restaurant is your mum''s friend. Okay, so you''re saying that you think you weren''t very good, and that''s because you were a bit slow and you got confused with some orders. But you said it''s common for people to do that on that first shift. That''s what she said to you, but...'
$ws.Cells.Item(20, 2).Value = 'But it''s probably just she''s just trying to make me feel better. She probably felt sorry for me because I''m obviously not very good at it, like I''m not very good at anything recently. I think she just felt sorry for me.'
$ws.Cells.Item(21, 1).Value = 'Okay, so I''ve just been noting down some of those thoughts that you''ve been having. You kind of say that you feel like you''re not very good, you''re slow, you''re confused with the orders, and you feel like you''re inconveniencing her and she feels sorry for you. I suppose I''m wondering, I mean, what makes you think that she feels sorry for you?'
$ws.Cells.Item(21, 2).Value = 'Like I said, she''s my mum''s friend, and she''s obviously feeling sorry for my mum as well because she''s got a daughter that''s not like anyone else of her age. And she''s just like all alone and she''s not getting involved. So I think she just feels sorry for me and she doesn''t want, like, to upset me. But yeah.'
$ws.Cells.Item(22, 1).Value = 'When you''re having these kind of thoughts, because that''s quite a lot of negative thoughts you''re describing, how do you feel?'
$ws.Cells.Item(22, 2).Value = 'Just feel empty, I don''t... I''m just calling as it is. I''m just calling it as reality. I feel like I''m being honest with myself by just voicing it and not trying to get around hiding it because I know, like, I can see it in myself.'
$ws.Cells.Item(23, 1).Value = 'So it sounds like you feel like those thoughts are very true then? I mean, if you were going to write them... So I''m going to pick out some of those thoughts you describe. You said ''I''m not very good'', and ''I''m inconveniencing them''. How much would you say you believe those on a scale out of 100, where 100 is like, I absolutely believe it and 0 is like, I don''t believe it at all?'
$ws.Cells.Item(23, 2).Value = 'I''m about 80, yeah.'
$ws.Cells.Item(24, 1).Value = 'So you feel quite sure that you definitely believe those, okay? And you said that when you''re thinking those thoughts, you start to feel quite empty.'
$ws.Cells.Item(24, 2).Value = 'I just feel distant and like I said, I''m not really myself, I just... It''s like I''m not really there, but when I am, I can''t even do a good job.'
$ws.Cells.Item(25, 1).Value = 'And are there any other emotions that you experience?'
$ws.Cells.Item(25, 2).Value = 'I feel drained. I''d say I feel like I''m not in touch with anything, just... it''s like everything''s going on around me, but I''m just not there. I''m just, it''s like I''m just watching.'
$ws.Cells.Item(26, 1).Value = 'Okay, so you feel empty, distant, you feel drained. And how, when you''re feeling these emotions, how do you feel in your body? Kind of what sensations do you notice?'
$ws.Cells.Item(26, 2).Value = 'I don''t know how to describe it. I get like, I have a bit of like, like, you know, when you kind of get like butterflies when you''re nervous and stuff, kind of get a bit of that.'
$ws.Cells.Item(27, 1).Value = 'You feel some butterflies? Yeah, did you have butterflies that night when you were at the restaurant?'
$ws.Cells.Item(27, 2).Value = 'Yeah, I was pretty nervous. Then when I got into it, they went away and like it distracted me for a bit. And so obviously I realized that I wasn''t doing a very good job.'
$ws.Cells.Item(28, 1).Value = 'Okay? And you, when you have those emotions, how does it affect your energy levels? Do you feel more energetic, less energetic, or about the same when you feel down?'
$ws.Cells.Item(28, 2).Value = 'Yeah, when you''re having those emotions that you were describing... I just feel drained, like didn''t even have like the energy to wake up in the morning sometimes, like I just felt pointless, like my body''s heavy and I don''t want to. Yeah.'
$ws.Cells.Item(29, 1).Value = 'So it sounds like you feel very tired. You''re struggling to wake up. You said it''s pointless. And what do you do then when you feel like that? So I guess you just talked about when you wake up in the morning and you feel heavy, and I suppose you''ve been maybe having some of these thoughts and those emotions potentially the night before, what do you then do? How do you respond to that?'
$ws.Cells.Item(29, 2).Value = 'I just have to get up I suppose. I don''t really have a choice because my mum, I live with my mom. I know it''s hard for her, so I don''t want her to be too worried. So I just... I just try to carry on as normal. It takes longer and I don''t want to, but I do, so I''ll go on with my routine.'
$ws.Cells.Item(30, 1).Value = 'So you''re a bit slower in getting up, but you are still getting up anyway. And what about your kind of, what do you think about when you''re doing that?'
$ws.Cells.Item(30, 2).Value = 'I''m just thinking that nobody else is like me my age, like nobody should be feeling like this. And I''m the only one that does.'
$ws.Cells.Item(31, 1).Value = 'So it sounds like you think a bit more about the problem and you kind of spend time thinking about...'
$ws.Cells.Item(31, 2).Value = 'How I hate it and I don''t want to feel like this, but I do, and nobody else I know does. So it makes me feel even more alone and like getting me even more down. And I just don''t want to disappoint anyone. I didn''t really want it.'
$ws.Cells.Item(32, 1).Value = 'Okay, so it sounds like you''re trying to toughen up on yourself a little bit, tell yourself to work harder, to be different.'
$ws.Cells.Item(32, 2).Value = 'Okay.'
$ws.Cells.Item(33, 1).Value = 'And when you start telling yourself that kind of thing, how does that make you feel? Does that make you feel better afterwards or worse?'
$ws.Cells.Item(33, 2).Value = 'Feels worse. Because I still can''t snap out of it, like I''m aware it''s happening, and I just, and I''m trying, but I just can''t. And then, and then it''s just another thing that I''m failing at, just another thing that I can''t really get.'

$ws = $wb.Worksheets.Item('enc_2')
$ws.Cells.Item(2, 1).Value = 'Hi Gabriella, it''s good to see you again. I appreciate you coming back for another session. Last time, you shared quite a bit about how you''ve been feeling since your breakup and the difficulties you''re experiencing with motivation and self-esteem. How have things been for you since we last talked?'
$ws.Cells.Item(2, 2).Value = 'Hi, um, thanks. I guess... things are kind of the same, honestly. Like, I still feel pretty stuck most days. I’ve been trying to get out of the house more, like doing that waitressing thing last week, but it didn’t really make me feel better, you know? Just more nervous and kind of overwhelmed. I don’t know if I’m making any progress or just going in circles.'
$ws.Cells.Item(3, 1).Value = 'It sounds like you''re putting in some effort by trying to get out and do new things, but you''re feeling that it''s not translating into feeling better, and maybe even making you feel more overwhelmed. When you think about the waitressing shift, what do you notice about your thoughts and feelings during or after that experience?'
$ws.Cells.Item(3, 2).Value = 'Yeah, I guess that’s pretty accurate. Um, during the shift, I was really nervous—like, my stomach was doing flips, and I kept thinking I was messing up. I kept telling myself I wasn’t good enough or that I was slow, and I felt kinda embarrassed when I didn’t get things right. Afterward, I just felt... more down, like I failed again, even though everyone said it’s normal for a first shift. I keep thinking maybe I shouldn’t even bother trying if I’m just gonna mess up anyway. It’s like, I want to feel better, but I just end up feeling worse.'
$ws.Cells.Item(4, 1).Value = 'It sounds like during that shift, your thoughts were quite negative—about your abilities and the way you’re perceived—and those thoughts led to feelings of embarrassment and disappointment afterward. You mentioned feeling like you’re just going in circles and that trying might be pointless because it leads to feeling worse. When you notice those thoughts, how much do you believe them on a scale from 0 to 100? For example, how strongly do you think "I''m just going to mess up" or "I''m not good enough"?'
$ws.Cells.Item(4, 2).Value = 'Um, I’d say... probably around 75 or 80. Like, I know part of me knows I might just be being hard on myself, but another part feels really convinced that I’m just not cut out for these things, or that I’m just going to keep messing up no matter what. It’s kind of frustrating because I want to believe I can get better, but those thoughts feel so real sometimes. And honestly, it’s kind of exhausting to keep fighting that voice in my head.'
$ws.Cells.Item(5, 1).Value = 'It’s understandable that those persistent, negative thoughts feel very convincing and exhausting. You’re caught between knowing they might not be entirely true and feeling overwhelmed by them anyway. When you think about trying to challenge these thoughts or tell yourself they might not be completely accurate, what kind of response do you usually have? Do you find it easy or hard to question or reframe those beliefs?'
$ws.Cells.Item(5, 2).Value = 'Honestly, it’s pretty hard. I mean, I try to tell myself, like, “Maybe I just need more practice,” or “Everyone makes mistakes,” but then I think, yeah, but what if I really am just not good enough? Or what if I’m just fooling myself? It’s like, I want to believe I can change, but my mind keeps pulling me back to those doubts. It’s kind of like, I want to be more hopeful, but I just don’t feel it yet.'
$ws.Cells.Item(6, 1).Value = 'That’s a very honest and insightful reflection. It sounds like you’re aware of the possibility of change and even trying to hold onto hope, but those deep-seated doubts keep pulling you back into a negative mindset. When you’re caught in that cycle of doubt and hope, what do you usually do? Do you find yourself ruminating, withdrawing, or perhaps trying to distract yourself?'
$ws.Cells.Item(6, 2).Value = 'Um, I think I mostly just ruminate, honestly. Like, I keep replaying everything in my head—what I said, what I did wrong, how I probably disappointed everyone, especially my mom. Sometimes I just want to distract myself, but I feel guilty about not doing anything productive or not trying hard enough to fix it. So, yeah, I guess I mostly just get stuck in my head, and that makes me feel even more tired and hopeless. It’s like I can’t really escape those thoughts, even if I want to.'
$ws.Cells.Item(7, 1).Value = 'It sounds like your mind tends to get stuck in a loop of these negative thoughts, especially about disappointing your mom and feeling tired and hopeless. That rumination seems to drain your energy and make it harder to see a way forward. When you notice yourself ruminating, what do you think might help, even just a little, to break that cycle? Or what have you tried that maybe offers some relief, even temporarily?'
$ws.Cells.Item(7, 2).Value = 'Honestly, I haven’t really figured out what works yet. Sometimes I try to distract myself by watching shows or scrolling on my phone, but then I just feel guilty because I’m not doing anything meaningful. Other times, I just sit there and think about everything, which makes it worse. I guess I haven’t found a good way to stop those thoughts from taking over. It’s like, I want to find something that can help me feel a little better, even if it’s just for a moment, but I don’t know what that is yet.'
$ws.Cells.Item(8, 1).Value = 'It makes sense that you''re seeking something to help you feel even briefly better, but haven''t quite found what works yet. That frustration is understandable. Sometimes, finding small, manageable activities that don’t feel overwhelming can be a helpful starting point. Have you noticed any activities or moments, even small ones, that have brought you a little relief or distraction in the past? Or perhaps things you’ve enjoyed before but haven''t done recently?'
$ws.Cells.Item(8, 2).Value = 'Hmm, I guess... I used to like drawing or writing in my journal, but I haven''t really done either in a while. When I was younger, I’d spend hours just sketching or writing stories, and it made me forget about everything for a bit. But now, I just feel like I don’t have the energy or motivation to do those things anymore. It’s like, even if I want to, I just can’t seem to get started or find the interest. I wonder if I should try again, but I’m not sure if it’ll make a difference.'
$ws.Cells.Item(9, 1).Value = 'Reconnecting with activities like drawing or journaling that once brought you some peace sounds like a meaningful idea. It’s understandable to feel hesitant, especially when motivation and energy are low. Sometimes, just allowing yourself to engage in these activities without pressure—maybe for a few minutes—can help gently break the cycle of rumination. How do you feel about giving yourself permission to try again, even if it’s just a small step?'
$ws.Cells.Item(9, 2).Value = 'I guess I could try, yeah. Like, maybe I don’t have to do it perfectly or spend hours on it. Just, you know, a little bit, just to see if it helps me feel a bit more... grounded or something. It’s just hard to believe it’ll actually make a difference, but I suppose I can give it a shot. I don’t want to keep feeling like this all the time, so maybe it’s worth trying.'
$ws.Cells.Item(10, 1).Value = 'That sounds like a very reasonable and compassionate approach—giving yourself permission to engage in those activities without the pressure of doing it perfectly. Even small steps can sometimes lead to noticeable changes over time. If you decide to try drawing or journaling, perhaps setting a timer for just a few minutes could make it feel more manageable. Would you be willing to share how it goes next time, or if you notice any small shifts in how you’re feeling?'
$ws.Cells.Item(10, 2).Value = 'Yeah, I think I could do that. I mean, it’s not like I have much else going on right now, so I might as well try. I guess I’ll start with just a few minutes, like you said, and see if it helps. It’s kind of embarrassing to even admit I’m trying this, but maybe it’s worth a shot. Thanks for suggesting it. I’ll let you know how it goes next time.'
$ws.Cells.Item(11, 1).Value = 'That’s a very brave and honest decision, Gabriella. Remember, taking even small steps like that is a positive move forward, and there’s no shame in trying. It’s all about finding what works for you, one step at a time. I look forward to hearing how it goes next time. For now, is there anything else you’d like to share or explore today?'
$ws.Cells.Item(11, 2).Value = 'Honestly, I’m not sure. I feel like I’ve been just going over the same stuff in my head, and I don’t know if there’s anything new I want to talk about. But maybe, I guess, I’m still struggling with feeling like I’m not good enough—like, even if I try those little things, I still doubt myself a lot. It’s kind of frustrating because I want to believe I can get better, but I just don’t know how to really believe it yet. Maybe that’s something I need to work on more.'
$ws.Cells.Item(12, 1).Value = 'It’s understandable to feel that way, Gabriella. Doubting yourself and feeling not good enough can be very persistent, especially when those thoughts have been part of your experience for a while. Recognizing that this is something you want to work on is an important insight. Building that belief in yourself often takes time and small successes, even if they feel minor at first. 
Would you be open to exploring some strategies or exercises that might help gradually strengthen your self-belief? Or perhaps discussing ways to challenge those doubts when they come up?'
$ws.Cells.Item(12, 2).Value = 'Yeah, I think I’d be willing to try that. I mean, I don’t really know how to challenge those doubts, though. It’s like, I hear them, and I just believe them most of the time. But maybe if I had some kind of way to, I don’t know, talk back to them or remind myself of other things, it could help. I just don’t really know where to start with that, or if it even works.'
$ws.Cells.Item(13, 1).Value = 'That’s a very honest and thoughtful reflection. Challenging those doubts can feel difficult, especially when they seem so convincing. A common approach in CBT is to identify these negative thoughts, examine the evidence for and against them, and then develop more balanced or realistic perspectives. It’s like having a conversation with that inner critic, where you gently question its assumptions and remind yourself of your strengths or positive qualities.
Would you like to work through an example of one of your doubts—say, the thought “I’m not good enough”—and see how we might challenge it together?'
$ws.Cells.Item(13, 2).Value = 'Yeah, I think that could be helpful. I mean, I guess I do have some moments where I realize I’m not as terrible as I think—I mean, I did get through the shift, even if I was nervous and slow. But then I just tell myself that it was just luck or that I was lucky they were patient with me. So, I don’t know, how do I really challenge that kind of thought? It feels like I’m just making excuses, but I don’t want to believe I’m just fooling myself either.'
$ws.Cells.Item(14, 1).Value = 'That’s a very insightful observation—you recognize your efforts and some moments of competence, but then dismiss them as luck or patience. Challenging that kind of thought involves looking at the evidence more objectively and considering alternative explanations. 
For example, when you think, “It was just luck,” you can ask yourself, “What evidence do I have that supports or contradicts that?” Have there been times when you managed tasks or handled situations well, even if you were nervous? Or perhaps moments where others acknowledged your effort or progress? 
Would you like to try framing some evidence for and against that thought?'
$ws.Cells.Item(14, 2).Value = 'Yeah, I suppose there have been times when I did okay, even if I was nervous. Like, I did remember some of the orders, and people didn’t seem too annoyed with me. And my mom said I was doing fine, even if I felt like I was messing up. But then I think, maybe she’s just being nice or trying to make me feel better. So, I guess I can see that there’s some truth in what you’re saying, but it’s hard to really believe it.'
$ws.Cells.Item(15, 1).Value = 'Exactly, Gabriella. You’re already doing a great job of examining the evidence—recognizing that you did some things right, and that others, like your mom, offered reassurance. The challenge is in balancing that evidence against the doubts that tell you otherwise. It’s normal for those doubts to feel very strong, especially when you’re feeling low. 
A helpful next step can be to create a more balanced thought—something like, “I may have been nervous and made some mistakes, but I also remembered some things and handled some situations okay.” Does that resonate with you? Would you like to try crafting a balanced thought together?'
$ws.Cells.Item(15, 2).Value = 'Yes, I think I could try that. It’s kind of strange, because part of me still feels like I’m just pretending or that I’m fooling myself, but I guess I can see that I did some things right. Maybe I can tell myself that even if I was nervous, I still managed to do some things okay, and that’s not nothing. I’ll try to keep that in mind next time those doubts come up. Thanks for helping me see that there’s a different way to look at it.'
$ws.Cells.Item(16, 1).Value = 'That’s a really thoughtful and balanced way to approach those doubts, Gabriella. Acknowledging both the difficulties and the successes can help you build a more realistic and compassionate view of yourself. Remember, challenging negative thoughts is a process, and it’s normal for doubts to linger. With practice, it can become easier to see your efforts and progress more clearly. 
Would you like to set a small goal for yourself before our next session, maybe related to practicing this kind of balanced thinking or engaging in one of your favorite activities again?'
$ws.Cells.Item(16, 2).Value = 'Yeah, I think I’d like to try that. Maybe I can practice challenging my doubts when they come up, just like we talked about, and also try to do a little drawing or writing in my journal. Even if it’s just for a few minutes, I think it might help me feel a bit more connected to myself. I’ll see how it goes, and I’ll let you know next time. Thanks for all the support today.'
$ws.Cells.Item(17, 1).Value = 'That sounds like a very practical and compassionate plan, Gabriella. Practicing challenging your doubts when they arise and reconnecting with your creative outlets are excellent steps toward feeling more grounded and confident. Remember, progress can be gradual, and every small effort counts. I look forward to hearing how it goes next time. 
Thank you for sharing honestly today. See you next time.'
$ws.Cells.Item(17, 2).Value = 'Thank you. I appreciate it. I’ll try my best. See you next time.'
$ws.Cells.Item(18, 1).Value = 'You''re very welcome, Gabriella. I admire your willingness to try and your thoughtful approach. Take care, and I look forward to hearing about your progress next time. See you then.'
$ws.Cells.Item(18, 2).Value = 'Goodbye.'

$ws = $wb.Worksheets.Item('enc_3')
$ws.Cells.Item(2, 1).Value = 'Hi Gabriella, it’s good to see you again. I appreciate you coming back for another session. Last time, you shared quite a bit about how you’ve been feeling since your breakup and the difficulties you’re experiencing with motivation and self-esteem. How have things been for you since we last talked?'
$ws.Cells.Item(2, 2).Value = 'Hi, thanks. Um, honestly, not much has changed. I still feel pretty stuck most days, you know? I tried going back to that waitressing thing last week, but honestly, it just made me more nervous and overwhelmed. I don’t know if I’m really making any progress or just going in circles, really. It’s kind of frustrating.'
$ws.Cells.Item(3, 1).Value = 'It sounds like despite your efforts to get out there and try new things, you''re feeling like it''s not helping much and perhaps even making you feel more overwhelmed. When you think about that waitressing shift, what thoughts or feelings come up during or after the experience?'
$ws.Cells.Item(3, 2).Value = 'Yeah, exactly. During the shift, I was, um, really nervous—like, my stomach was doing flips, and I kept thinking I was messing up. I kept telling myself I wasn’t good enough or that I was slow, and I felt kind of embarrassed when I didn’t get things right. Afterward, I just felt… more down, like I failed again, even though everyone said it’s normal for a first shift. I keep thinking maybe I shouldn’t even bother trying if I’m just gonna mess up anyway. It’s like, I want to feel better, but I just end up feeling worse.'
$ws.Cells.Item(4, 1).Value = 'It sounds like during the shift, your thoughts were quite negative—about your abilities and the way you''re perceived—and those thoughts led to feelings of embarrassment and disappointment afterward. You mentioned feeling like you''re just going in circles and that trying might be pointless because it leads to feeling worse. When you notice those thoughts, how much do you believe them on a scale from 0 to 100? For example, how strongly do you think "I''m just going to mess up" or "I''m not good enough"?'
$ws.Cells.Item(4, 2).Value = 'Um, I’d say… probably around 75 or 80. Like, I know part of me knows I might just be being hard on myself, but another part feels really convinced that I’m just not cut out for these things, or that I’m just going to keep messing up no matter what. It’s kind of frustrating because I want to believe I can get better, but those thoughts feel so real sometimes. And honestly, it’s kind of exhausting to keep fighting that voice in my head.'
$ws.Cells.Item(5, 1).Value = 'That’s very understandable. Those persistent, negative thoughts can feel very convincing and exhausting. You’re caught between knowing they might not be entirely true and feeling overwhelmed by them anyway. When you think about trying to challenge these thoughts or tell yourself they might not be completely accurate, what kind of response do you usually have? Do you find it easy or hard to question or reframe those beliefs?'
$ws.Cells.Item(5, 2).Value = 'Honestly, it’s pretty hard. I mean, I try to tell myself, like, “Maybe I just need more practice,” or “Everyone makes mistakes,” but then I think, yeah, but what if I really am just not good enough? Or what if I’m just fooling myself? It’s like, I want to believe I can change, but my mind keeps pulling me back to those doubts. It’s kind of like, I want to be more hopeful, but I just don’t feel it yet.'
$ws.Cells.Item(6, 1).Value = 'That’s a very honest and insightful reflection. It sounds like you’re aware of the possibility of change and even trying to hold onto hope, but those deep-seated doubts keep pulling you back into a negative mindset. When you’re caught in that cycle of doubt and hope, what do you usually do? Do you find yourself ruminating, withdrawing, or perhaps trying to distract yourself?'
$ws.Cells.Item(6, 2).Value = 'Yeah, I think I mostly just ruminate, honestly. Like, I keep replaying everything in my head—what I said, what I did wrong, how I probably disappointed everyone, especially my mom. Sometimes I just want to distract myself, but I feel guilty about not doing anything productive or not trying hard enough to fix it. So, yeah, I guess I mostly just get stuck in my head, and that makes me feel even more tired and hopeless. It’s like, I can’t really escape those thoughts, even if I want to.'
$ws.Cells.Item(7, 1).Value = 'It sounds like your mind tends to get stuck in a loop of these negative thoughts, especially about disappointing your mom and feeling tired and hopeless. That rumination seems to drain your energy and make it harder to see a way forward. When you notice yourself ruminating, what do you think might help, even just a little, to break that cycle? Or what have you tried that maybe offers some relief, even temporarily?'
$ws.Cells.Item(7, 2).Value = 'Honestly, I haven’t really figured out what works yet. Sometimes I try to distract myself by watching shows or scrolling on my phone, but then I just feel guilty because I’m not doing anything meaningful. Other times, I just sit there and think about everything, which makes it worse. I guess I haven’t found a good way to stop those thoughts from taking over. It’s like, I want to find something that can help me feel a little better, even if it’s just for a moment, but I don’t know what that is yet.'
$ws.Cells.Item(8, 1).Value = 'It makes sense that you''re seeking something to help you feel even briefly better, but haven''t quite found what works yet. That frustration is understandable. Sometimes, finding small, manageable activities that don’t feel overwhelming can be a helpful starting point. Have you noticed any activities or moments, even small ones, that have brought you a little relief or distraction in the past? Or perhaps things you’ve enjoyed before but haven''t done recently?'
$ws.Cells.Item(8, 2).Value = 'Hmm, I guess... I used to like drawing or writing in my journal, but I haven''t really done either in a while. When I was younger, I’d spend hours just sketching or writing stories, and it made me forget about everything for a bit. But now, I just feel like I don’t have the energy or motivation to do those things anymore. It’s like, even if I want to, I just can’t seem to get started or find the interest. I wonder if I should try again, but I’m not sure if it’ll make a difference.'
$ws.Cells.Item(9, 1).Value = 'Reconnecting with activities like drawing or journaling that once brought you some peace sounds like a meaningful idea. It’s understandable to feel hesitant, especially when motivation and energy are low. Sometimes, just allowing yourself to engage in these activities without pressure—maybe for a few minutes—can help gently break the cycle of rumination. How do you feel about giving yourself permission to try again, even if it’s just a small step?'
$ws.Cells.Item(9, 2).Value = 'I guess I could try, yeah. Like, maybe I don’t have to do it perfectly or spend hours on it. Just, you know, a little bit, just to see if it helps me feel a bit more... grounded or something. It’s just hard to believe it’ll actually make a difference, but I suppose I can give it a shot. I don’t want to keep feeling like this all the time, so maybe it’s worth trying.'
$ws.Cells.Item(10, 1).Value = 'That sounds like a very reasonable and compassionate approach—giving yourself permission to engage in those activities without the pressure of doing it perfectly. Even small steps can sometimes lead to noticeable changes over time. If you decide to try drawing or journaling, perhaps setting a timer for just a few minutes could make it feel more manageable. Would you be willing to share how it goes next time, or if you notice any small shifts in how you’re feeling?'
$ws.Cells.Item(10, 2).Value = 'Yeah, I think I could do that. I mean, it’s not like I have much else going on right now, so I might as well try. I guess I’ll start with just a few minutes, like you said, and see if it helps. It’s kind of embarrassing to even admit I’m trying this, but maybe it’s worth a shot. Thanks for suggesting it. I’ll let you know how it goes next time.'
$ws.Cells.Item(11, 1).Value = 'That’s a very brave and honest decision, Gabriella. Remember, taking even small steps like that is a positive move forward, and there’s no shame in trying. It’s all about finding what works for you, one step at a time. I look forward to hearing how it goes next time. For now, is there anything else you’d like to share or explore today?'
$ws.Cells.Item(11, 2).Value = 'I guess, um, I don’t know. I feel like I’ve just been going over the same stuff in my head, and I’m not really sure if there’s anything new I want to talk about. But I guess I’m still struggling with feeling like I’m not good enough—like, even if I do those little things, I still doubt myself a lot. It’s kind of frustrating because I want to believe I can get better, but I just don’t know how to really believe it yet. Maybe that’s something I need to work on more.'
$ws.Cells.Item(12, 1).Value = 'It’s completely understandable to feel that way, Gabriella. Doubting yourself and feeling not good enough can be very persistent, especially when those thoughts have been part of your experience for a while. Recognizing that this is something you want to work on is an important insight. Building that belief in yourself often takes time and small successes, even if they feel minor at first.
Would you be open to exploring some strategies or exercises that might help gradually strengthen your self-belief? Or perhaps discussing ways to challenge those doubts when they come up?'
$ws.Cells.Item(12, 2).Value = 'Yeah, I think I’d be willing to try that. I mean, I don’t really know how to challenge those doubts, though. It’s like, I hear them, and I just believe them most of the time. But maybe if I had some kind of way to, I don’t know, talk back to them or remind myself of other things, it could help. I just don’t really know where to start with that, or if it even works.'
$ws.Cells.Item(13, 1).Value = 'That’s a very honest and thoughtful reflection. Challenging those doubts can feel difficult, especially when they seem so convincing. A common approach in CBT is to identify these negative thoughts, examine the evidence for and against them, and then develop more balanced or realistic perspectives. It’s like having a conversation with that inner critic, where you gently question its assumptions and remind yourself of your strengths or positive qualities.
Would you like to work through an example of one of your doubts—say, the thought “I’m not good enough”—and see how we might challenge it together?'
$ws.Cells.Item(13, 2).Value = 'Yeah, I think that could be helpful. I mean, I guess I do have some moments where I realize I’m not as terrible as I think—I mean, I did get through the shift, even if I was nervous and slow. But then I just tell myself that it was just luck or that I was lucky they were patient with me. So, I don’t know, how do I really challenge that kind of thought? It feels like, I’m just making excuses, but I don’t want to believe I’m just fooling myself either.'
$ws.Cells.Item(14, 1).Value = 'That’s a very insightful observation—you recognize your efforts and some moments of competence, but then dismiss them as luck or patience. Challenging that kind of thought involves looking at the evidence more objectively and considering alternative explanations. 
For example, when you think, “It was just luck,” you can ask yourself, “What evidence do I have that supports or contradicts that?” Have there been times when you managed tasks or handled situations well, even if you were nervous? Or perhaps moments where others acknowledged your effort or progress? 
Would you like to try framing some evidence for and against that thought?'
$ws.Cells.Item(14, 2).Value = 'Yeah, I suppose there have been times when I did okay, even if I was nervous. Like, I did remember some of the orders, and people didn’t seem too annoyed with me. And my mom said I was doing fine, even if I felt like I was messing up. But then I think, maybe she’s just being nice or trying to make me feel better. So, I guess I can see that there’s some truth in what you’re saying, but it’s hard to really believe it.'
$ws.Cells.Item(15, 1).Value = 'Exactly, Gabriella. You’re already doing a great job of examining the evidence—recognizing that you did some things right, and that others, like your mom, offered reassurance. The challenge is in balancing that evidence against the doubts that tell you otherwise. It’s normal for those doubts to feel very strong, especially when you’re feeling low.
A helpful next step can be to create a more balanced thought—something like, “I may have been nervous and made some mistakes, but I also remembered some things and handled some situations okay.” Does that resonate with you? Would you like to try crafting a balanced thought together?'
$ws.Cells.Item(15, 2).Value = 'Yes, I think I could try that. It’s kind of strange, because part of me still feels like I’m just pretending or that I’m fooling myself, but I guess I can see that I did some things right. Maybe I can tell myself that even if I was nervous, I still managed to do some things okay, and that’s not nothing. I’ll try to keep that in mind next time those doubts come up. Thanks for helping me see that there’s a different way to look at it.'
$ws.Cells.Item(16, 1).Value = 'That’s a very thoughtful and balanced way to approach those doubts, Gabriella. Acknowledging both the difficulties and the successes can help you build a more realistic and compassionate view of yourself. Remember, challenging negative thoughts is a process, and it’s normal for doubts to linger. With practice, it can become easier to see your efforts and progress more clearly.
Would you like to set a small goal for yourself before our next session, maybe related to practicing this kind of balanced thinking or engaging in one of your favorite activities again?'
$ws.Cells.Item(16, 2).Value = 'Yeah, I think I’d like to try that. Maybe I can practice challenging my doubts when they come up, just like we talked about, and also try to do a little drawing or writing in my journal. Even if it’s just for a few minutes, I think it might help me feel a bit more connected to myself. I’ll see how it goes, and I’ll let you know next time. Thanks for all the support today.'
$ws.Cells.Item(17, 1).Value = 'That sounds like a very practical and compassionate plan, Gabriella. Practicing challenging your doubts when they arise and reconnecting with your creative outlets are excellent steps toward feeling more grounded and confident. Remember, progress can be gradual, and every small effort counts. I look forward to hearing how it goes next time.
Thank you for sharing honestly today. See you next time.'
$ws.Cells.Item(17, 2).Value = 'Thank you. I appreciate it. I’ll try my best. See you next time.'
$ws.Cells.Item(18, 1).Value = 'You''re very welcome, Gabriella. I admire your willingness to try and your thoughtful approach. Take care, and I look forward to hearing about your progress next time. See you then.'
$ws.Cells.Item(18, 2).Value = 'Goodbye.'

$wb.Worksheets.Item("enc_4").Delete()

